$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "January 2021 " (in the Tyler Technologies entry, "January 2021 -
# Present") becomes "January 2022 ". In the canonical XML this is expressed
# as a run split: "January 202" stays in the original run, and the final
# digit "2" is joined with the run that used to hold just the trailing
# space, becoming a new run containing "2 ".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("January 2021") | Out-Null
# Range covering the trailing "1" of the year plus the space that follows it.
$tail = $d.Range($rng.End - 1, $rng.End + 1)
# Toggle formatting on/off around the write so the engine is forced to open
# a fresh run here instead of silently re-merging it back with its neighbor.
$tail.Font.Bold = 1
$tail.Text = "2 "
$tail.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: the separate "– " and "Remote" runs (after "Tyler Technologies ")
# collapse into a single run "– Remote" (same text, same formatting, just
# no longer split across two runs). Because the visible text is not
# changing, first stamp in placeholder text (forcing a real content
# change/new run) and then restore the real text, otherwise a same-value
# write is a no-op that leaves the original two runs untouched.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("– Remote") | Out-Null
$placeholder = $d.Range($rng2.Start, $rng2.End)
$placeholder.Font.Bold = 1
$placeholder.Text = "XXXXXXXX"
$merge = $d.Range($rng2.Start, $rng2.Start + 8)
$merge.Text = "– Remote"
$merge.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 3: the second "January 2021" (end of the "November 2019 – January
# 2021" line) also becomes "January 2022", again expressed as a run split:
# "January 202" keeps the original run, and a brand-new trailing run holds
# just the final "2".
# ---------------------------------------------------------------------------
$rng3 = $d.Content
while ($rng3.Find.Execute("January 2021", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
  $lastEnd = $rng3.End
  $rng3.Collapse(0) | Out-Null
}
$tail3 = $d.Range($lastEnd - 1, $lastEnd)
$tail3.Font.Bold = 1
$tail3.Text = "2"
$tail3.Font.Bold = 0
